# Update tutorial workbook:
# - insert new "immigration" sheet (with citizenship axis) between "deaths" and "pop_births_deaths"
# - update population figures for France (2014 / 2015) across all sheets that contain them
# - update the aggregated "pop_narrow_format" totals for France (2014 / 2015)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "immigration" worksheet right after "deaths"
# ---------------------------------------------------------------------------
$sheetDeaths = $wb.Worksheets.Item("deaths")
$immigration = $wb.Worksheets.Add($null, $sheetDeaths)
$immigration.Name = "immigration"

# header row
$immigration.Cells.Item(1, 1).Value = "country"
$immigration.Cells.Item(1, 2).Value = "citizenship"
$immigration.Cells.Item(1, 3).Value = "gender\time"
$immigration.Cells.Item(1, 4).Value = 2013
$immigration.Cells.Item(1, 5).Value = 2014
$immigration.Cells.Item(1, 6).Value = 2015

$immigrationData = @(
  ,@("Belgium","Belgium","Male",8822,10512,11378)
  ,@("Belgium","Belgium","Female",5727,6301,6486)
  ,@("Belgium","Luxembourg","Male",102,117,105)
  ,@("Belgium","Luxembourg","Female",117,123,114)
  ,@("Belgium","Netherlands","Male",4185,4222,4183)
  ,@("Belgium","Netherlands","Female",3737,3844,3942)
  ,@("Luxembourg","Belgium","Male",896,937,880)
  ,@("Luxembourg","Belgium","Female",574,655,622)
  ,@("Luxembourg","Luxembourg","Male",694,722,660)
  ,@("Luxembourg","Luxembourg","Female",607,586,535)
  ,@("Luxembourg","Netherlands","Male",160,165,147)
  ,@("Luxembourg","Netherlands","Female",92,97,85)
  ,@("Netherlands","Belgium","Male",1063,1141,1113)
  ,@("Netherlands","Belgium","Female",980,1071,1181)
  ,@("Netherlands","Luxembourg","Male",23,43,59)
  ,@("Netherlands","Luxembourg","Female",24,34,46)
  ,@("Netherlands","Netherlands","Male",19374,20037,21119)
  ,@("Netherlands","Netherlands","Female",16945,17411,18084)
)

$r = 2
foreach ($row in $immigrationData) {
  $c = 1
  foreach ($val in $row) {
    $immigration.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update France population figures (2014 / 2015) on every sheet that has
#    the full "pop" table (France Male / France Female rows 4 & 5)
# ---------------------------------------------------------------------------
$popSheetNames = @("pop", "pop_births_deaths", "pop_missing_axis_name")
foreach ($name in $popSheetNames) {
  $sheet = $wb.Worksheets.Item($name)
  # France / Male -> row 4 ; France / Female -> row 5
  $sheet.Cells.Item(4, 4).Value = 32045129
  $sheet.Cells.Item(4, 5).Value = 32174258
  $sheet.Cells.Item(5, 4).Value = 34120851
  $sheet.Cells.Item(5, 5).Value = 34283895
}

# pop_missing_values has no "France / Male" row, so France / Female is row 4
$sheetMissing = $wb.Worksheets.Item("pop_missing_values")
$sheetMissing.Cells.Item(4, 4).Value = 34120851
$sheetMissing.Cells.Item(4, 5).Value = 34283895

# ---------------------------------------------------------------------------
# 3. Update the France totals (2014 / 2015) in pop_narrow_format
# ---------------------------------------------------------------------------
$sheetNarrow = $wb.Worksheets.Item("pop_narrow_format")
$sheetNarrow.Cells.Item(6, 3).Value = 66165980
$sheetNarrow.Cells.Item(7, 3).Value = 66458153

# keep "pop_narrow_format" as the selected / active sheet (tab), as in the source file
$sheetNarrow.Activate()
